$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New column K header: date 2021-05-26 (serial 44342), formatted like the other date cells in row 4
$ws.Range("J4").Copy()
$ws.Range("K4").PasteSpecial(-4122)
$ws.Range("K4").Value = 44342

# New column K body (rows 5-15): check marks, formatted like column J's check marks
for ($r = 5; $r -le 15; $r++) {
    $srcCell = $ws.Cells.Item($r, 10)
    $dstCell = $ws.Cells.Item($r, 11)
    $srcCell.Copy()
    $dstCell.PasteSpecial(-4122)
    $dstCell.Value = [char]0x2713
}

$excel.CutCopyMode = 0

# Update selection to K15, as in the target workbook
$ws.Range("K15").Select()
